$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked" and "is_enabled" template columns (F:G), which held
# the is_locked_lbl / is_enabled_lbl shared-string templates. Deleting the
# entire columns shifts the remaining columns (order_by, rem) left and lets
# Excel drop the now-unreferenced shared strings automatically.
$ws.Range("F1:G1").EntireColumn.Delete()
